$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff: each entry sets the target
# text value for a cell on the "cryptos" worksheet (price/volume/hour columns).
# A leading apostrophe forces Excel to store the value as literal text (so
# numeric-looking strings like "300.20", "-6.50%" or "4" are preserved exactly,
# matching the original inline-string cell type), and resetting the cell Style
# back to "Normal" afterwards clears the quote-prefix formatting flag that
# entering a text-forced value would otherwise leave behind.
$updates = @(
    @{ Cell = "D2"; Value = '300.20' },
    @{ Cell = "E2"; Value = '-6.50%' },
    @{ Cell = "G2"; Value = '4' },
    @{ Cell = "D3"; Value = '35.06' },
    @{ Cell = "E3"; Value = '-3.04%' },
    @{ Cell = "G3"; Value = '4' },
    @{ Cell = "D4"; Value = '4.985' },
    @{ Cell = "E4"; Value = '-3.30%' },
    @{ Cell = "G4"; Value = '4' },
    @{ Cell = "D5"; Value = '0.07920' },
    @{ Cell = "E5"; Value = '-2.09%' },
    @{ Cell = "G5"; Value = '4' },
    @{ Cell = "D6"; Value = '1.925' },
    @{ Cell = "E6"; Value = '-10.84%' },
    @{ Cell = "G6"; Value = '4' },
    @{ Cell = "D7"; Value = '7.738' },
    @{ Cell = "E7"; Value = '-4.09%' },
    @{ Cell = "G7"; Value = '4' },
    @{ Cell = "D8"; Value = '4.011' },
    @{ Cell = "E8"; Value = '-2.86%' },
    @{ Cell = "G8"; Value = '4' },
    @{ Cell = "E9"; Value = '4.48%' },
    @{ Cell = "G9"; Value = '4' },
    @{ Cell = "D10"; Value = '0.9241' },
    @{ Cell = "E10"; Value = '-0.45%' },
    @{ Cell = "G10"; Value = '4' },
    @{ Cell = "D11"; Value = '0.1137' },
    @{ Cell = "E11"; Value = '12.59%' },
    @{ Cell = "G11"; Value = '4' },
    @{ Cell = "D12"; Value = '0.1825' },
    @{ Cell = "E12"; Value = '-3.17%' },
    @{ Cell = "G12"; Value = '4' },
    @{ Cell = "D13"; Value = '0.09175' },
    @{ Cell = "E13"; Value = '-0.90%' },
    @{ Cell = "G13"; Value = '4' },
    @{ Cell = "D14"; Value = '0.03528' },
    @{ Cell = "E14"; Value = '-1.19%' },
    @{ Cell = "G14"; Value = '4' },
    @{ Cell = "D15"; Value = '0.09890' },
    @{ Cell = "E15"; Value = '-0.45%' },
    @{ Cell = "G15"; Value = '4' },
    @{ Cell = "D16"; Value = '0.001398' },
    @{ Cell = "E16"; Value = '-2.65%' },
    @{ Cell = "G16"; Value = '4' },
    @{ Cell = "D17"; Value = '0.005790' },
    @{ Cell = "E17"; Value = '2.48%' },
    @{ Cell = "G17"; Value = '4' },
    @{ Cell = "D18"; Value = '3.508' },
    @{ Cell = "E18"; Value = '1.54%' },
    @{ Cell = "G18"; Value = '4' },
    @{ Cell = "D19"; Value = '0.3442' },
    @{ Cell = "E19"; Value = '2.05%' },
    @{ Cell = "G19"; Value = '4' },
    @{ Cell = "E20"; Value = '-1.56%' },
    @{ Cell = "G20"; Value = '4' },
    @{ Cell = "D21"; Value = '5.071' },
    @{ Cell = "E21"; Value = '-0.32%' },
    @{ Cell = "G21"; Value = '4' },
    @{ Cell = "E22"; Value = '8.87%' },
    @{ Cell = "G22"; Value = '4' },
    @{ Cell = "D23"; Value = '0.04492' },
    @{ Cell = "E23"; Value = '-2.48%' },
    @{ Cell = "G23"; Value = '4' },
    @{ Cell = "E24"; Value = '-2.20%' },
    @{ Cell = "G24"; Value = '4' },
    @{ Cell = "G25"; Value = '4' },
    @{ Cell = "D26"; Value = '0.0001250' },
    @{ Cell = "E26"; Value = '-3.89%' },
    @{ Cell = "G26"; Value = '4' },
    @{ Cell = "E27"; Value = '-6.77%' },
    @{ Cell = "G27"; Value = '4' },
    @{ Cell = "G28"; Value = '4' },
    @{ Cell = "G29"; Value = '4' },
    @{ Cell = "G30"; Value = '4' },
    @{ Cell = "G31"; Value = '4' },
    @{ Cell = "G32"; Value = '4' },
    @{ Cell = "G33"; Value = '4' },
    @{ Cell = "G34"; Value = '4' },
    @{ Cell = "G35"; Value = '4' },
    @{ Cell = "G36"; Value = '4' },
    @{ Cell = "G37"; Value = '4' },
    @{ Cell = "G38"; Value = '4' },
    @{ Cell = "D39"; Value = '0.01881' },
    @{ Cell = "E39"; Value = '-6.56%' },
    @{ Cell = "G39"; Value = '4' },
    @{ Cell = "D40"; Value = '0.04691' },
    @{ Cell = "E40"; Value = '-5.98%' },
    @{ Cell = "G40"; Value = '4' },
    @{ Cell = "D41"; Value = '0.007599' },
    @{ Cell = "E41"; Value = '-2.31%' },
    @{ Cell = "G41"; Value = '4' },
    @{ Cell = "D42"; Value = '0.009554' },
    @{ Cell = "E42"; Value = '22.07%' },
    @{ Cell = "G42"; Value = '4' },
    @{ Cell = "E43"; Value = '-5.78%' },
    @{ Cell = "G43"; Value = '4' },
    @{ Cell = "D44"; Value = '0.002120' },
    @{ Cell = "E44"; Value = '1.83%' },
    @{ Cell = "G44"; Value = '4' },
    @{ Cell = "D45"; Value = '0.01110' },
    @{ Cell = "E45"; Value = '-7.82%' },
    @{ Cell = "G45"; Value = '4' },
    @{ Cell = "D46"; Value = '0.00006014' },
    @{ Cell = "E46"; Value = '-6.18%' },
    @{ Cell = "G46"; Value = '4' },
    @{ Cell = "E47"; Value = '0.04%' },
    @{ Cell = "G47"; Value = '4' },
    @{ Cell = "E48"; Value = '123.04%' },
    @{ Cell = "G48"; Value = '4' },
    @{ Cell = "E49"; Value = '-31.31%' },
    @{ Cell = "G49"; Value = '4' },
    @{ Cell = "D50"; Value = '0.00002100' },
    @{ Cell = "E50"; Value = '0.04%' },
    @{ Cell = "G50"; Value = '4' },
    @{ Cell = "D51"; Value = '0.0002000' },
    @{ Cell = "E51"; Value = '0.04%' },
    @{ Cell = "G51"; Value = '4' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.Value = "'" + $u.Value
    $range.Style = "Normal"
}
